# Discord.docx edit: turn the bare invite URL into a real hyperlink
# (new invite code) and split the paragraph so the "_GoBack" bookmark
# ends up on its own paragraph, exactly as the target OOXML shows.

$d = $word.ActiveDocument

# --- 1. Swap the URL text and push the bookmark into a new paragraph ---
# The whole document body is currently one paragraph:
#   "https://discord.gg/XKZjq" + bookmarkStart/_GoBack + bookmarkEnd
# Replacing just the visible text (not the trailing paragraph mark) with
# the new URL followed by a paragraph break ("`r") turns that into two
# paragraphs: the URL text in paragraph 1, and the (still-present)
# bookmark left behind in the new paragraph 2 - matching the diff.
$oldText = "https://discord.gg/XKZjq"
$newText = "https://discord.gg/YfKJuya"

$urlRange = $d.Range(0, $oldText.Length)
$urlRange.Text = $newText + "`r"

# --- 2. Materialize the built-in "Hyperlink" character style ---
# Word stamps a concrete w:style definition for "Hyperlink" (blue,
# underlined) the first time it is used; recreate that definition so the
# run we are about to convert can reference it via rStyle.
$style = $d.Styles.Add("Hyperlink", 2)
$style.BaseStyle = "DefaultParagraphFont"
$style.Priority = 99
$style.UnhideWhenUsed = $true
$style.Font.Underline = 1
$style.Font.TextColor.ObjectThemeColor = 10

# --- 3. Turn the new URL text into an actual hyperlink field ---
$linkRange = $d.Range(0, $newText.Length)
$d.Hyperlinks.Add($linkRange, $newText) | Out-Null

Write-Output "Converted invite link to hyperlink and split paragraph."
